# Horarios actualizados Línea 141 - 556
# Updates the "scraped" timestamp and arrival data for the three sheets
# (LP1912, LP1912-215, 6203-6173) of the schedule workbook.

$wb = $excel.ActiveWorkbook

$nuevaHora = "04:40:49"

# ---------------------------------------------------------------------------
# Sheet "LP1912": full table refresh, 3 new rows appended (9 -> 12 rows)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: " + $nuevaHora
$ws1.Range("A3").Value = "Total filas: 12"

$filas1 = @(
    @($nuevaHora, "04:47", "81_EL PELIGRO",   7, "LP1912"),
    @($nuevaHora, "04:53", "11_ETCHEVERRY",  13, "LP1912"),
    @($nuevaHora, "05:16", "17_ROMERO",      36, "LP1912"),
    @($nuevaHora, "05:22", "23_HERNANDEZ",   42, "LP1912"),
    @($nuevaHora, "05:44", "14_ABASTO",      64, "LP1912"),
    @($nuevaHora, "05:46", "17_ROMERO",      66, "LP1912"),
    @($nuevaHora, "06:00", "16_SANTA ANA",   80, "LP1912"),
    @($nuevaHora, "06:09", "10_OLMOS",       89, "LP1912"),
    @($nuevaHora, "06:15", "215A_EL PATO",   95, "LP1912"),
    @($nuevaHora, "06:30", "23_HERNANDEZ",  110, "LP1912"),
    @($nuevaHora, "06:34", "11_ETCHEVERRY", 114, "LP1912"),
    @($nuevaHora, "06:38", "17X38_ROMERO",  118, "LP1912")
)

$fila = 6
foreach ($registro in $filas1) {
    $ws1.Cells.Item($fila, 1).Value = $registro[0]
    $ws1.Cells.Item($fila, 2).Value = $registro[1]
    $ws1.Cells.Item($fila, 3).Value = $registro[2]
    $ws1.Cells.Item($fila, 4).Value = $registro[3]
    $ws1.Cells.Item($fila, 5).Value = $registro[4]
    $fila = $fila + 1
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": only the timestamp and the single data row refresh
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: " + $nuevaHora

$ws2.Cells.Item(6, 1).Value = $nuevaHora
$ws2.Cells.Item(6, 2).Value = "06:15"
$ws2.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(6, 4).Value = 95
$ws2.Cells.Item(6, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": only the timestamp refreshes, no data rows
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: " + $nuevaHora
